# Regenerate merged AHB files
#  - rename header columns: "*_old" -> "*_FV2410", "*_new" -> "*_FV2504"
#  - freeze header row (row 1)
#  - wrap the data range in a native Excel Table (ListObject) with AutoFilter

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row -------------------------------------------------
# Columns A:J  -> "<name>_old" becomes "<name>_FV2410"
# Column  K    -> "diff" (unchanged)
# Columns L:U  -> "<name>_new" becomes "<name>_FV2504"
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2410")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2504")
}

# --- 2. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the data range into a native Table (ListObject) --------------
$dataRange = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
